$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "43.008.98"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -5.33%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.224.19"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -6.27%  "

$ws.Range("E4").Value = "  -0.04%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "313.95"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "100.35"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -7.34%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.584"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -7.39%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -8.55%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.77"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -9.91%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "54.67"
$cell.Style = "Normal"

$ws.Range("E12").Value = "  -10.42%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "7.60"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -10.54%  "

$ws.Range("E14").Value = "  -1.34%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.561.61"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -6.39%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.854"
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.17"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -7.62%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.222.33"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -7.21%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "42.898.25"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -5.59%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.13"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("E21").Value = "  -9.74%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.48"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -10.76%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "65.54"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -10.70%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -12.29%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "237.03"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -8.78%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.10"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -12.36%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  +1.76%  "

$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("E30").Value = "  -11.08%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.38"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -11.66%  "

$ws.Range("E32").Value = "  -8.83%  "

$ws.Range("E33").Value = "  -11.80%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "34.27"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -7.93%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "154.22"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -7.67%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.76"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -7.17%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.04"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.09%  "

$ws.Range("E38").Value = "  -7.00%  "

$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("E40").Value = "  -6.79%  "

$ws.Range("E41").Value = "  -11.96%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.66"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -8.62%  "

$ws.Range("E43").Value = "  -9.52%  "

$ws.Range("E44").Value = "  -0.07%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "12.57"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.797.12"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.29%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "85.83"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -13.57%  "

$ws.Range("E48").Value = "  -10.78%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "76.20"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -9.53%  "

$ws.Range("E50").Value = "  -8.52%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "58.98"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -15.77%  "

